$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '55.267.97'
$ws.Range("E2").Value = '  -3.60%  '
$ws.Range("D3").Value = '2.952.48'
$ws.Range("E3").Value = '  -5.97%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '487.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.69%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.14'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.47%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '2.950.95'
$ws.Range("E8").Value = '  -6.09%  '
$ws.Range("E9").Value = '  -6.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.04'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.22%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.100'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.96%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.347'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -9.02%  '
$ws.Range("E13").Value = '  +0.37%  '
$ws.Range("D14").Value = '3.456.44'
$ws.Range("E14").Value = '  -5.78%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '24.58'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.22%  '
$ws.Range("D16").Value = '55.177.07'
$ws.Range("E16").Value = '  -3.75%  '
$ws.Range("D17").Value = '2.947.54'
$ws.Range("E17").Value = '  -6.01%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000139'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -6.69%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.54'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.99%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.00'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -6.29%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.42'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '315.12'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -7.82%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.461'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -8.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '60.03'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -12.23%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.997'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.160'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.68%  '
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("D29").Value = '0.0₃0838'
$ws.Range("E29").Value = '  -9.58%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.45'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.09%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.50'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.71%  '
$ws.Range("E32").Value = '  -1.71%  '
$ws.Range("E33").Value = '  -8.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.27'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -10.37%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '148.95'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.34%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.35'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -9.30%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.28'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.96%  '
$ws.Range("B38").Value = 'Aptos'
$ws.Range("C38").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.64'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -8.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '23.20'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -9.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0643'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.45%  '
$ws.Range("D41").Value = '2.980.52'
$ws.Range("E41").Value = '  -5.86%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '35.91'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -10.65%  '
$ws.Range("E44").Value = '  -6.36%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.628'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -8.68%  '
$ws.Range("E46").Value = '  -6.58%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.52'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -8.50%  '
$ws.Range("D48").Value = '2.121.46'
$ws.Range("E48").Value = '  -4.81%  '
$ws.Range("E49").Value = '  +0.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.94'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.84%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.54'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -9.16%  '
